$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.510.10"
$ws.Range("E2").Value = "  +0.20%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.821.03"
$ws.Range("E3").Value = "  -0.26%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.15"
$ws.Range("E5").Value = "  -0.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5101"
$ws.Range("E7").Value = "  -5.57%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3949"
$ws.Range("E8").Value = "  -2.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08325"
$ws.Range("E9").Value = "  +8.70%  "

$ws.Range("E10").Value = "  -0.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.66"
$ws.Range("E11").Value = "  -0.56%  "

$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.13"
$ws.Range("E12").Value = "  +1.03%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.323"
$ws.Range("E13").Value = "  -0.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.000"
$ws.Range("E14").Value = "  -0.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.539"
$ws.Range("E15").Value = "  -1.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.817.59"
$ws.Range("E16").Value = "  -0.58%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001141"
$ws.Range("E17").Value = "  +5.89%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.67"
$ws.Range("E18").Value = "  +3.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06653"
$ws.Range("E19").Value = "  +0.91%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.78"
$ws.Range("E20").Value = "  +0.62%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9995"
$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.099"
$ws.Range("E22").Value = "  +0.47%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.545.04"
$ws.Range("E23").Value = "  +0.31%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.47"
$ws.Range("E24").Value = "  +3.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.266"
$ws.Range("E25").Value = "  +1.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.28"
$ws.Range("E26").Value = "  +2.66%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.95"
$ws.Range("E27").Value = "  -0.85%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.025.73"
$ws.Range("E28").Value = "  -0.61%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.412"
$ws.Range("E29").Value = "  -2.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.80"
$ws.Range("E30").Value = "  +1.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.110"
$ws.Range("E31").Value = "  -1.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1095"
$ws.Range("E32").Value = "  -2.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.786"
$ws.Range("E33").Value = "  +1.93%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.653"
$ws.Range("E34").Value = "  +0.34%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07075"
$ws.Range("E35").Value = "  -4.12%  "

$ws.Range("E36").Value = "  -0.86%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02340"
$ws.Range("E37").Value = "  +0.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.228"
$ws.Range("E38").Value = "  +0.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.860"
$ws.Range("E39").Value = "  +0.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6302"
$ws.Range("E40").Value = "  +0.51%  "

$ws.Range("E41").Value = "  -0.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.178"
$ws.Range("E42").Value = "  -0.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9994"
$ws.Range("E43").Value = "  -0.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.399"
$ws.Range("E44").Value = "  +0.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.48"
$ws.Range("E45").Value = "  -0.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5928"
$ws.Range("E46").Value = "  +1.36%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.729"
$ws.Range("E47").Value = "  +0.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.32"
$ws.Range("E48").Value = "  +0.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.987"
$ws.Range("E49").Value = "  -0.89%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.187"
$ws.Range("E50").Value = "  -1.29%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06893"
$ws.Range("E51").Value = "  +0.09%  "
